$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value would otherwise be auto-interpreted by Excel as a
# number (e.g. "1.001", "0.9974") need NumberFormat forced to Text ("@") first
# so they are stored as strings, matching the source data which are all text.
$textCells = @(
    "D4",
    "D5",
    "D6",
    "D7",
    "D8",
    "D9",
    "D10",
    "D11",
    "D12",
    "D13",
    "D14",
    "D15",
    "D17",
    "D18",
    "D19",
    "D20",
    "D22",
    "D23",
    "D25",
    "D26",
    "D27",
    "D28",
    "D29",
    "D31",
    "D32",
    "D33",
    "D34",
    "D35",
    "D36",
    "D37",
    "D38",
    "D39",
    "D40",
    "D41",
    "D42",
    "D43",
    "D44",
    "D45",
    "D46",
    "D47",
    "D48",
    "D49",
    "D50",
    "D51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "24.288.23"
$ws.Range("E2").Value = "  +11.64%  "

# Row 3
$ws.Range("D3").Value = "1.683.67"
$ws.Range("E3").Value = "  +6.78%  "

# Row 4
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.67%  "

# Row 5
$ws.Range("D5").Value = "310.13"
$ws.Range("E5").Value = "  +9.08%  "

# Row 6
$ws.Range("D6").Value = "0.9974"
$ws.Range("E6").Value = "  +2.97%  "

# Row 7
$ws.Range("D7").Value = "0.3747"
$ws.Range("E7").Value = "  +1.88%  "

# Row 8
$ws.Range("D8").Value = "0.3456"
$ws.Range("E8").Value = "  +5.86%  "

# Row 9
$ws.Range("D9").Value = "47.79"
$ws.Range("E9").Value = "  +16.64%  "

# Row 10
$ws.Range("D10").Value = "1.199"
$ws.Range("E10").Value = "  +5.82%  "

# Row 11
$ws.Range("D11").Value = "0.07340"
$ws.Range("E11").Value = "  +4.24%  "

# Row 12
$ws.Range("D12").Value = "0.9975"
$ws.Range("E12").Value = "  +0.92%  "

# Row 13
$ws.Range("D13").Value = "20.62"
$ws.Range("E13").Value = "  +2.62%  "

# Row 14
$ws.Range("D14").Value = "6.135"
$ws.Range("E14").Value = "  +5.89%  "

# Row 15
$ws.Range("D15").Value = "6.800"
$ws.Range("E15").Value = "  +4.98%  "

# Row 16
$ws.Range("D16").Value = "1.679.42"
$ws.Range("E16").Value = "  +7.71%  "

# Row 17
$ws.Range("D17").Value = "0.00001117"
$ws.Range("E17").Value = "  +4.60%  "

# Row 18
$ws.Range("D18").Value = "0.9970"
$ws.Range("E18").Value = "  +3.03%  "

# Row 19
$ws.Range("D19").Value = "0.06734"
$ws.Range("E19").Value = "  +9.29%  "

# Row 20
$ws.Range("D20").Value = "82.24"
$ws.Range("E20").Value = "  +11.58%  "

# Row 21
$ws.Range("E21").Value = "  +4.15%  "

# Row 22
$ws.Range("D22").Value = "6.150"
$ws.Range("E22").Value = "  +5.75%  "

# Row 23
$ws.Range("D23").Value = "12.07"
$ws.Range("E23").Value = "  +4.15%  "

# Row 24
$ws.Range("D24").Value = "24.232.03"
$ws.Range("E24").Value = "  +11.86%  "

# Row 25
$ws.Range("D25").Value = "2.409"
$ws.Range("E25").Value = "  +3.92%  "

# Row 26
$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").Value = "2.693"
$ws.Range("E26").Value = "  +11.72%  "

# Row 27
$ws.Range("B27").Value = "LEO"
$ws.Range("C27").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D27").Value = "3.363"
$ws.Range("E27").Value = "  -9.00%  "

# Row 28
$ws.Range("D28").Value = "152.12"
$ws.Range("E28").Value = "  +2.51%  "

# Row 29
$ws.Range("D29").Value = "19.66"
$ws.Range("E29").Value = "  +8.42%  "

# Row 30
$ws.Range("D30").Value = "1.863.26"
$ws.Range("E30").Value = "  +7.84%  "

# Row 31
$ws.Range("D31").Value = "127.34"
$ws.Range("E31").Value = "  +6.66%  "

# Row 32
$ws.Range("D32").Value = "6.551"
$ws.Range("E32").Value = "  +22.67%  "

# Row 33
$ws.Range("D33").Value = "4.061"
$ws.Range("E33").Value = "  +0.86%  "

# Row 34
$ws.Range("D34").Value = "1.000"
$ws.Range("E34").Value = "  +11.53%  "

# Row 35
$ws.Range("D35").Value = "1.784"
$ws.Range("E35").Value = "  +14.88%  "

# Row 36
$ws.Range("D36").Value = "0.08512"
$ws.Range("E36").Value = "  +4.55%  "

# Row 37
$ws.Range("D37").Value = "12.64"
$ws.Range("E37").Value = "  +10.20%  "

# Row 38
$ws.Range("D38").Value = "0.06500"
$ws.Range("E38").Value = "  +8.45%  "

# Row 39
$ws.Range("D39").Value = "5.413"
$ws.Range("E39").Value = "  +6.77%  "

# Row 40
$ws.Range("D40").Value = "8.968"
$ws.Range("E40").Value = "  +9.71%  "

# Row 41
$ws.Range("D41").Value = "0.02362"
$ws.Range("E41").Value = "  +9.28%  "

# Row 42
$ws.Range("D42").Value = "1.289"
$ws.Range("E42").Value = "  +4.48%  "

# Row 43
$ws.Range("D43").Value = "0.2140"
$ws.Range("E43").Value = "  +6.60%  "

# Row 44
$ws.Range("D44").Value = "0.6226"
$ws.Range("E44").Value = "  +8.16%  "

# Row 45
$ws.Range("D45").Value = "0.9972"
$ws.Range("E45").Value = "  +3.09%  "

# Row 46
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "13.29"
$ws.Range("E46").Value = "  +4.42%  "

# Row 47
$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").Value = "3.820"
$ws.Range("E47").Value = "  +6.21%  "

# Row 48
$ws.Range("D48").Value = "0.5979"
$ws.Range("E48").Value = "  +5.95%  "

# Row 49
$ws.Range("D49").Value = "127.95"
$ws.Range("E49").Value = "  +3.46%  "

# Row 50
$ws.Range("D50").Value = "2.038"
$ws.Range("E50").Value = "  +6.66%  "

# Row 51
$ws.Range("D51").Value = "0.07175"
$ws.Range("E51").Value = "  +6.74%  "
